# Parametrize city on certificates
# Replace the static "Latacunga, " literal before the current_date MERGEFIELD
# with a new MERGEFIELD for =consultation.branch_office.city, keeping the
# trailing ", " as a separate literal run.

$d = $word.ActiveDocument

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("Latacunga, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Latacunga, ' run to replace"
}

$insertAt = $rng.Start

# Remove the old literal text, collapsing the range to its start.
$rng.Text = ""

# Build the replacement field markup: a complex MERGEFIELD for the branch
# office city, immediately followed by the literal ", " that used to trail
# "Latacunga".
$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD =consultation.branch_office.city \* MERGEFORMAT </w:instrText></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>«=consultation.branch_office.city»</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>
</w:p>
'@

$insertRng = $d.Range($insertAt, $insertAt)
$insertRng.InsertXML($newXml)
